# Updated cryptos list on Wed Mar 27 21:47:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.198.67"
$ws.Range("E2").Value = "'  -1.25%  "
$ws.Range("D3").Value = "'3.516.07"
$ws.Range("E3").Value = "'  -2.12%  "
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'573.80"
$ws.Range("E5").Value = "'  -0.98%  "
$ws.Range("D6").Value = "'184.72"
$ws.Range("D7").Value = "'3.507.58"
$ws.Range("E7").Value = "'  -2.28%  "
$ws.Range("E8").Value = "'  -2.91%  "
$ws.Range("E9").Value = "'  +0.09%  "
$ws.Range("E10").Value = "'  +1.21%  "
$ws.Range("D11").Value = "'0.654"
$ws.Range("E11").Value = "'  -1.61%  "
$ws.Range("D12").Value = "'54.35"
$ws.Range("E12").Value = "'  -3.16%  "
$ws.Range("E13").Value = "'  -1.68%  "
$ws.Range("E14").Value = "'  -2.06%  "
$ws.Range("D15").Value = "'4.080.07"
$ws.Range("D16").Value = "'19.47"
$ws.Range("E16").Value = "'  -2.67%  "
$ws.Range("D17").Value = "'69.146.83"
$ws.Range("E17").Value = "'  -1.40%  "
$ws.Range("D18").Value = "'3.514.58"
$ws.Range("E18").Value = "'  -2.11%  "
$ws.Range("D19").Value = "'12.35"
$ws.Range("E19").Value = "'  -2.72%  "
$ws.Range("D20").Value = "'0.120"
$ws.Range("E20").Value = "'  -1.17%  "
$ws.Range("D21").Value = "'542.18"
$ws.Range("E21").Value = "'  +13.47%  "
$ws.Range("E22").Value = "'  -2.85%  "
$ws.Range("E23").Value = "'  -5.97%  "
$ws.Range("D24").Value = "'5.01"
$ws.Range("E24").Value = "'  -0.95%  "
$ws.Range("D25").Value = "'4.43"
$ws.Range("D26").Value = "'94.43"
$ws.Range("E26").Value = "'  -1.18%  "
$ws.Range("D27").Value = "'11.13"
$ws.Range("E27").Value = "'  +0.19%  "
$ws.Range("E28").Value = "'  -2.98%  "
$ws.Range("D29").Value = "'9.13"
$ws.Range("E29").Value = "'  -3.59%  "
$ws.Range("D30").Value = "'31.84"
$ws.Range("E30").Value = "'  -1.71%  "
$ws.Range("D31").Value = "'7.26"
$ws.Range("E31").Value = "'  -5.41%  "
$ws.Range("D32").Value = "'12.66"
$ws.Range("E32").Value = "'  +3.38%  "
$ws.Range("D33").Value = "'64.69"
$ws.Range("E33").Value = "'  -2.67%  "
$ws.Range("E34").Value = "'  -4.36%  "
$ws.Range("D35").Value = "'559.03"
$ws.Range("E35").Value = "'  -4.74%  "
$ws.Range("D36").Value = "'3.13"
$ws.Range("E36").Value = "'  +8.69%  "
$ws.Range("D37").Value = "'38.13"
$ws.Range("E37").Value = "'  -2.62%  "
$ws.Range("D38").Value = "'0.403"
$ws.Range("E38").Value = "'  +1.52%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "'  -0.13%  "
$ws.Range("D40").Value = "'0.0₃0765"
$ws.Range("E40").Value = "'  -5.23%  "
$ws.Range("B41").Value = "'Stacks"
$ws.Range("C41").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.38"
$ws.Range("E41").Value = "'  -2.43%  "
$ws.Range("B42").Value = "'dogwifhat"
$ws.Range("C42").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.09"
$ws.Range("E42").Value = "'  -4.66%  "
$ws.Range("E43").Value = "'  -3.81%  "
$ws.Range("D44").Value = "'3.280.94"
$ws.Range("E44").Value = "'  +1.37%  "
$ws.Range("E45").Value = "'  -3.03%  "
$ws.Range("D46").Value = "'0.0446"
$ws.Range("E46").Value = "'  +0.20%  "
$ws.Range("D47").Value = "'3.49"
$ws.Range("E47").Value = "'  +3.94%  "
$ws.Range("E48").Value = "'  -2.60%  "
$ws.Range("D49").Value = "'8.91"
$ws.Range("E49").Value = "'  -5.68%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "'  -0.26%  "
$ws.Range("D51").Value = "'137.60"
$ws.Range("E51").Value = "'  +1.84%  "
